# The commit inserts one new daily price record for "Apio" (Vega Modelo de
# Temuco) right above the existing row 215, pushing every subsequent row
# down by one (old row 215 -> new row 216, ..., old row 338 -> new row 339).
# All columns of the shifted rows keep their original values; only the new
# row 215 carries fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 215 - this shifts rows 215..338 down to
# 216..339 and copies formatting (incl. the date style on column D) from
# the row above, matching how the row was originally added in Excel.
$ws.Rows(215).Insert()

# Populate the newly inserted row with the new record's data. Columns that
# are constant across this sheet (mercado/region/category/unit/etc.) reuse
# the same values as the neighboring rows.
$ws.Cells.Item(215, 1).Value2  = 10
$ws.Cells.Item(215, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(215, 3).Value2  = "La Araucanía"
$ws.Cells.Item(215, 4).Value2  = 44777
$ws.Cells.Item(215, 5).Value2  = 9
$ws.Cells.Item(215, 6).Value2  = 100112017
$ws.Cells.Item(215, 7).Value2  = "Apio"
$ws.Cells.Item(215, 8).Value2  = "Americana (o)"
$ws.Cells.Item(215, 9).Value2  = "Primera"
$ws.Cells.Item(215, 10).Value2 = 110
$ws.Cells.Item(215, 11).Value2 = 12000
$ws.Cells.Item(215, 12).Value2 = 12000
$ws.Cells.Item(215, 13).Value2 = 12000
$ws.Cells.Item(215, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(215, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(215, 16).Value2 = 2000
$ws.Cells.Item(215, 17).Value2 = 6
$ws.Cells.Item(215, 18).Value2 = "Hortaliza"
